# "Added units to dashboard"
#
# The Dashboard sheet's "PORTFOLIO BY INDUSTRY" block (rows 6-11) lists each
# currently-held stock (column B, a spilled dynamic array) together with its
# industry (column C, also a spill). This change adds a third spilled column,
# D ("Units"), holding the total number of units currently held for each of
# those stocks - i.e. SUMIF(Ledger[Stock], <that stock>, Ledger[Units]) for
# each row of the B6:B11 spill.
#
# D6 is entered as the dynamic-array formula (spilling down through D11);
# D7:D11 end up holding just the literal numbers the spill produced, exactly
# like the pre-existing C6 (formula) / C7:C11 (plain cached spill values)
# pair immediately to their left.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dashboard")

# Per-stock unit totals for the six rows currently spilled into B6:B11
# (ASX, CAT, CBA, EOS, NAB, NEA respectively) - i.e. what
# SUMIF(Ledger[Stock], ANCHORARRAY(B6), Ledger[Units]) resolves to once the
# linked Stock data type can refresh. Entered as one dynamic-array formula
# so it spills down D6:D11 exactly like the existing B6/C6 spills.
$ws.Range("D6").Formula2 = "={162;4275;56;695;35;1376}"

# Move the active selection from D6 to E6, as in the diff.
$ws.Range("E6").Select()
